$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

$newCaseQuery = @"
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
WHERE demo.breed IN ['Shih Tzu'] 
MATCH (c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS ``Case ID`` ,
        coalesce(s.clinical_study_designation, '') AS ``Study Code`` ,
        coalesce(s.clinical_study_type, '') AS  ``Study Type``,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS ``Stage of Disease`` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS ``Neutered Status``,
        coalesce(demo.weight, '') AS ``Weight (kg)``,
        coalesce(diag.best_response, '') AS ``Response to Treatment``
"@

$ws1.Range("B2").Value = $newCaseQuery
$ws1.Range("B2").Select()
